$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Overview sheet: status text changed from "Ready for handoff" to
# "Handed back: in sync with en-US" for both zh-cn and de-de columns.
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# Columns widened to fit the longer status text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ----------------------------------------------------------------------
# zh-cn sheet: status updated, handback datetime refreshed, error detail
# cleared (file is now in sync, so the "not latest" warning goes away).
# ----------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-11-03 19:17:53"
$wsZhCn.Range("K3").Value = "2016-11-03 19:17:53"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

# ----------------------------------------------------------------------
# de-de sheet: status updated, handback datetime refreshed, error detail
# cleared.
# ----------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-11-03 19:18:12"
$wsDeDe.Range("K3").Value = "2016-11-03 19:18:12"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334

Write-Output "Generated handback report"
